$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "CountPairsWithGivenSum"
$ws.Range("A15").Value = "Count pairs with given sum"

$ws.Range("A15").Select()
